# Updates cryptos list values (Price and Volume(1h) columns) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.940.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4695"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3663"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07161"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9276"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07708"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.882.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.289"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.405"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008625"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.963.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.028"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.920"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.021"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.880"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08859"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.224"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7463"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.784"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.478"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.086"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01940"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05194"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.911"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1520"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4696"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.609"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06047"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8927"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.29%  "

